$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "hip knee and joint",
    "hip knee joints",
    "hip leggings",
    "hip model",
    "hip pack",
    "hip pack black",
    "hip pack for women",
    "hip pack for women hiking",
    "hip pack medical",
    "hip pack women",
    "hip pain",
    "hip pain compression",
    "hip pain right side",
    "hip pain women",
    "hip pain yoga",
    "hip pocket",
    "hip pocket for women",
    "hip replacement gifts",
    "hip ring",
    "hip slimmer",
    "hip slimmer for women",
    "hip stabilizer and groin brace",
    "hip supplements women",
    "hip support",
    "hip support for women",
    "hip support pants",
    "hip support shorts women",
    "hip supports",
    "hip warmers for women",
    "hip womens clothing",
    "hips and knees",
    "hipster crop top",
    "hipster pants for women",
    "hipster womens clothing",
    "hit tennis 3",
    "hit the spot game",
    "hobble dress long",
    "hole tights",
    "hologlamour leggings",
    "holy clothing for women",
    "holy clothing plus",
    "holy cow",
    "holy leggings",
    "holy leggings for women",
    "home gift guide",
    "hoop frisbee",
    "hoop slips for under dresses",
    "hoops vs mid",
    "horse ridding clothes",
    "hose and tights",
    "hose for women",
    "hot chilly women",
    "hot chillys leggings",
    "hot hoops replacement balls",
    "hot leggings for women",
    "hot n nasty",
    "hot run",
    "hot spots soccer",
    "hot weather clothing for women",
    "hot weather gear",
    "hot weather pants",
    "how women rise",
    "hug pot",
    "hyperextension back",
    "i day delivery",
    "i do crew tank",
    "i give thanks",
    "i hip",
    "i in team",
    "i love butter",
    "i love colors tights",
    "i love target",
    "i over soccer ball",
    "i put pants on today",
    "i red thailand",
    "i rise black motion",
    "i roast 2",
    "i target pro",
    "i walk 2",
    "ibs support",
    "ideal stretch hamstring",
    "iliac crest",
    "iliotibial band",
    "iliotibial band syndrome",
    "impact shorts women",
    "in the company of cheerful ladies",
    "incline exercise ball",
    "increased circulation",
    "inhaler bronchitis",
    "inhaler wick",
    "inhaler wicks",
    "inhalers for bronchitis",
    "inner game of basketball",
    "inner outer thigh machine",
    "inner tennis playing the game",
    "inner thigh slimmer",
    "insulated bicycle pants",
    "insulated clothing for women",
    "insulated exercise pants",
    "insulated hiking pants women"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
